$d = $word.ActiveDocument

# 1) Paragraph "O problema da dificuldade em controlar o estoque afeta os técnicos..."
#    Shorten the "afeta" clause and rewrite the "devido" clause.
$found1 = $d.Content.Find.Execute(
    " os técnicos, que muitas vezes precisam desmarcar serviços pois não possuem os materiais necessários ao atendimento ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " os técnicos ", 2)

$found2 = $d.Content.Find.Execute(
    " a não possuírem previamente informações sobre o estoque de peças.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " não atenderem ", 2)

# 2) The paragraph right after that one was empty; it now gets the continuation text.
$pFill = $d.Paragraphs.Item(3)
$pFill.Range.Text = "as ordens de serviços e consequentemente não receberem a quantidade de comissão que poderiam, por conta do estoque desatualizado."

# 3) Insert a brand-new empty paragraph right after it (a spacer before the next
#    "O problema ..." block), matching the new blank line in the document.
$pFill.Range.InsertParagraphAfter() | Out-Null

# 4) Paragraph "O problema do negócio não ter uma grande visibilidade..."
#    Rewrite the "devido" clause.
$found3 = $d.Content.Find.Execute(
    " a empresa não conseguir, por falta de divulgação, alcançar as metas de lucro desejadas.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " não conseguir alcançar as metas de lucro desejadas.", 2)

# 5) Make the page orientation explicit (portrait) on the section page size.
$d.PageSetup.Orientation = 0

Write-Output "replace1=$found1 replace2=$found2 replace3=$found3 paragraphs=$($d.Paragraphs.Count)"
